# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.704.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.785.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.73%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.783.63'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.75%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E10").Value = '  -0.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.452'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.04%  '

$ws.Range("E13").Value = '  +4.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.417.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.777.86'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.645.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.19%  '

$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("E21").Value = '  -6.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '468.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.718'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.73%  '

$ws.Range("E24").Value = '  -7.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.49%  '

$ws.Range("E26").Value = '  -1.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.43%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.934.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.62%  '

$ws.Range("E33").Value = '  -2.71%  '

$ws.Range("E34").Value = '  -3.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.747.25'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.97%  '

$ws.Range("E38").Value = '  -0.27%  '

$ws.Range("E39").Value = '  -1.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.138'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.312'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.69%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.65'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.37%  '

$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '395.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.40%  '

$ws.Range("E49").Value = '  -8.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.66'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.57%  '
